$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text in the source data (values such as
# "64.711.24" or "6.18" are not valid numbers). Excel auto-converts a
# plain .Value assignment that looks numeric into a real number, which
# would lose the original text formatting (e.g. "6.20" -> 6.2). To avoid
# that, force the cell to Text format while writing the value, then
# restore the cell's normal style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.865.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.33%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.428.90'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.13%  '

$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.10%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.428.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.556'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.97%  '

$ws.Range("E10").Value = '  +1.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.121'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.99%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.424'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.017.88'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.29%  '

$ws.Range("E14").Value = '  +0.23%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.46%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000174'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.855.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.405.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.88'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.93%  '

$ws.Range("E25").Value = '  -2.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000117'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.41%  '

$ws.Range("E28").Value = '  +0.19%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.02'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.03'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.40%  '

$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.68%  '

$ws.Range("E36").Value = '  -5.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.872'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.97%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.85'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0726'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.778.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.82'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("E44").Value = '  +0.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.41%  '

$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.15%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0306'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.98%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("E50").Value = '  -0.64%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.98%  '
